# Update automàtic: dades i banners [2026-02-04 16:38]
# Updates the DATA_EXTRACCIO (column E) timestamps for rows 2-36.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$timestamps = @{
    2  = "2026-02-04 16:37:04"
    3  = "2026-02-04 16:37:06"
    4  = "2026-02-04 16:37:09"
    5  = "2026-02-04 16:37:12"
    6  = "2026-02-04 16:37:15"
    7  = "2026-02-04 16:37:17"
    8  = "2026-02-04 16:37:20"
    9  = "2026-02-04 16:37:23"
    10 = "2026-02-04 16:37:26"
    11 = "2026-02-04 16:37:28"
    12 = "2026-02-04 16:37:31"
    13 = "2026-02-04 16:37:34"
    14 = "2026-02-04 16:37:36"
    15 = "2026-02-04 16:37:39"
    16 = "2026-02-04 16:37:42"
    17 = "2026-02-04 16:37:45"
    18 = "2026-02-04 16:37:47"
    19 = "2026-02-04 16:37:50"
    20 = "2026-02-04 16:37:53"
    21 = "2026-02-04 16:37:56"
    22 = "2026-02-04 16:37:58"
    23 = "2026-02-04 16:38:01"
    24 = "2026-02-04 16:38:04"
    25 = "2026-02-04 16:38:07"
    26 = "2026-02-04 16:38:09"
    27 = "2026-02-04 16:38:12"
    28 = "2026-02-04 16:38:15"
    29 = "2026-02-04 16:38:18"
    30 = "2026-02-04 16:38:21"
    31 = "2026-02-04 16:38:23"
    32 = "2026-02-04 16:38:26"
    33 = "2026-02-04 16:38:28"
    34 = "2026-02-04 16:38:31"
    35 = "2026-02-04 16:38:33"
    36 = "2026-02-04 16:38:36"
}

foreach ($row in $timestamps.Keys) {
    $ws.Cells.Item($row, 5).Value = $timestamps[$row]
}
